# Update "想去人数" (F column) counts across the workbook's sheets to
# reflect the latest scrape (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 33
$ws1.Range("F3").Value = 494
$ws1.Range("F7").Value = 1317
$ws1.Range("F10").Value = 1356
$ws1.Range("F12").Value = 1100
$ws1.Range("F13").Value = 191
$ws1.Range("F15").Value = 263
$ws1.Range("F16").Value = 1707
$ws1.Range("F17").Value = 632
$ws1.Range("F19").Value = 305
$ws1.Range("F20").Value = 3341
$ws1.Range("F21").Value = 29
$ws1.Range("F22").Value = 417
$ws1.Range("F23").Value = 941
$ws1.Range("F24").Value = 1223
$ws1.Range("F25").Value = 1897
$ws1.Range("F27").Value = 1678
$ws1.Range("F31").Value = 672
$ws1.Range("F33").Value = 10
$ws1.Range("F34").Value = 1952
$ws1.Range("F35").Value = 916
$ws1.Range("F36").Value = 1967
$ws1.Range("F37").Value = 216
$ws1.Range("F38").Value = 470
$ws1.Range("F39").Value = 130
$ws1.Range("F41").Value = 352
$ws1.Range("F42").Value = 921
$ws1.Range("F43").Value = 821
$ws1.Range("F44").Value = 1063
$ws1.Range("F45").Value = 156
$ws1.Range("F48").Value = 237

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F13").Value = 821

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 494
$ws4.Range("F8").Value = 1317
$ws4.Range("F11").Value = 1356
$ws4.Range("F13").Value = 1100
$ws4.Range("F14").Value = 191
$ws4.Range("F16").Value = 263
$ws4.Range("F17").Value = 1707
$ws4.Range("F18").Value = 632
$ws4.Range("F20").Value = 305
$ws4.Range("F21").Value = 3341
$ws4.Range("F22").Value = 29
$ws4.Range("F23").Value = 417
$ws4.Range("F25").Value = 1223
$ws4.Range("F27").Value = 1678
$ws4.Range("F30").Value = 821
$ws4.Range("F34").Value = 1952
$ws4.Range("F35").Value = 916
$ws4.Range("F36").Value = 1967
$ws4.Range("F37").Value = 470
$ws4.Range("F38").Value = 130
$ws4.Range("F40").Value = 921
$ws4.Range("F41").Value = 821
$ws4.Range("F42").Value = 1063
$ws4.Range("F43").Value = 156
$ws4.Range("F47").Value = 237
